$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-20 Sunday" "2024-10-21 Monday"

Replace-Text "60×95=5700" "24×84=2016"
Replace-Text "70×22=1540" "65×22=1430"
Replace-Text "42×90=3780" "82×23=1886"
Replace-Text "98×89=8722" "79×71=5609"
Replace-Text "30×49=1470" "52×27=1404"
Replace-Text "85×82=6970" "29×54=1566"
Replace-Text "49×22=1078" "61×36=2196"
Replace-Text "11×82=902" "52×99=5148"
Replace-Text "44×25=1100" "72×61=4392"
Replace-Text "57×79=4503" "41×77=3157"
Replace-Text "35×54=1890" "49×53=2597"
Replace-Text "12×27=324" "78×98=7644"
Replace-Text "93×59=5487" "41×34=1394"
Replace-Text "14×54=756" "82×17=1394"
Replace-Text "56×19=1064" "25×29=725"
Replace-Text "34×30=1020" "60×16=960"
Replace-Text "56×78=4368" "18×55=990"
Replace-Text "35×24=840" "91×87=7917"
Replace-Text "13×92=1196" "65×74=4810"
Replace-Text "55×26=1430" "59×90=5310"
Replace-Text "24×20=480" "20×69=1380"
Replace-Text "44×92=4048" "11×24=264"
Replace-Text "39×84=3276" "35×40=1400"
Replace-Text "75×39=2925" "92×56=5152"
Replace-Text "82×92=7544" "32×55=1760"
